# Daily attendance processing - 2026-02-01 02:18:58
# Normalizes the "Recorded By" (column G) value so that "System" is listed
# before the recorder's email address.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $current = $cell.Value()
    if ($current -eq $oldValue) {
        $cell.Value = $newValue
    }
}
